# "Add files via upload" — re-upload of app/ByCoach.xlsx with a handful of
# "Started" (Yes/No) values corrected for a few players, plus the window
# left scrolled a little further up the list before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# --- Fix the "Started" (column C) flags that were toggled for these rows ---
# Row 10  (Tramon Mark, Booz)            Yes -> No
# Row 12  (Billy Richmond III, Booz)     No  -> Yes
# Row 45  (Boogie Fland, Hal)            Yes -> No
# Row 46  (AJ Storr, Hal)                No  -> Yes
# Row 66  (Josh Hubbard, Ron)            No  -> Yes
# Row 68  (Rylan Griffen, Ron)           Yes -> No
# Row 77  (Mike Sharavjamts, Tar)        Yes -> No
# Row 84  (Tahaad Pettiford, Tar)        No  -> Yes
$ws.Range("C10").Value = "No"
$ws.Range("C12").Value = "Yes"
$ws.Range("C45").Value = "No"
$ws.Range("C46").Value = "Yes"
$ws.Range("C66").Value = "Yes"
$ws.Range("C68").Value = "No"
$ws.Range("C77").Value = "No"
$ws.Range("C84").Value = "Yes"

# --- Scroll the frozen (header-row) view down so row 59 is the first
#     visible row under the frozen header, matching where the sheet was
#     left when it was re-saved ---
$win = $excel.ActiveWindow
$win.ScrollRow = 59
$win.ScrollColumn = 1

# Restore the original active selection
$ws.Range("C85").Select()
